$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff/Handback Datetime for the
# fe7f28da-cc97-4606-83ad-7a7b8a3ad11d row (row 5)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-17 03:25:44"
$wsZhCn.Range("G5").Value = "2016-02-17 03:26:27"

# de-de sheet: update Correspond Handoff/Handback Datetime for the
# fe7f28da-cc97-4606-83ad-7a7b8a3ad11d row (row 5)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-17 03:25:54"
$wsDeDe.Range("G5").Value = "2016-02-17 03:26:44"
